$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "https://villarta.com.br/elevadores-e-escadas-rolantes-villarta/lista-de-empresas-de-elevadores-em-sp/"
$ws.Range("B2").Value = "protecaodedados@villarta.com.br;"

# Add new rows 4-9
$ws.Range("A4").Value = "https://spelevadores.com.br/"
$ws.Range("B4").Value = "contato@spelevadores.com.br;"

$ws.Range("A5").Value = "https://www.basselevadores.com.br/elevadores-sao-paulo-sp.php"
$ws.Range("B5").Value = "contato@basselevadores.com.br;"

$ws.Range("A6").Value = "https://elevadoresoiwa.com.br/"
$ws.Range("B6").Value = "comercial@elevadoresoiwa.com.br;oiwa@elevadoresoiwa.com.br;"

$ws.Range("A7").Value = "https://www.otis.com/pt/br"
$ws.Range("B7").Value = "Navigati_ouvidoria@otis.com;ouvidoria@otis.com;Navigati_imprensa@otis.com;cac@otis.com;imprensa@otis.com;"

$ws.Range("A8").Value = "https://crel.com.br/"
$ws.Range("B8").Value = "bruno@crel.com.br;"

$ws.Range("A9").Value = "https://www.primac.com.br/manutencao-de-elevadores-sp.php"
$ws.Range("B9").Value = "comercial@primac.com.br;"

# Apply the same style as existing data rows (style index 2, which corresponds to A2/B2) to the new rows
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B9").PasteSpecial(-4122)  # xlPasteFormats
